$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 133
$ws1.Range("F3").Value = 209
$ws1.Range("F4").Value = 3569
$ws1.Range("F5").Value = 370
$ws1.Range("G5").Value = 58
$ws1.Range("F6").Value = 20
$ws1.Range("F7").Value = 431

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 5

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 133
$ws4.Range("F3").Value = 209
$ws4.Range("F4").Value = 3569
$ws4.Range("F5").Value = 370
$ws4.Range("G5").Value = 58
$ws4.Range("F7").Value = 5
$ws4.Range("F8").Value = 20
$ws4.Range("F9").Value = 431
